# global_mean_data sheet: switch from a "year-per-row" table (TB, PC, dS0,
# dSF, K, IP, I0, IF as columns) to a "stat-per-row" table (2021..2024, mean
# as columns), using K = Strike# instead of the prior K values, and
# regenerated std/mean figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Propagate the header style (bold font + thin border + centered/top
#        alignment, currently style index 1 on B1:I1 and A2:A6) onto the three
#        extra label rows we are about to add (A7:A9), *before* touching any
#        other content, by copying formats only from an already-styled cell.
$ws.Range("B1").Copy()
$ws.Range("A7:A9").PasteSpecial(-4122)  # xlPasteFormats

# --- 2. Wipe the old numeric data body (B2:I6) - it will be fully replaced
#        with regenerated values in the new column layout.
$ws.Range("B2:I6").ClearContents()

# --- 3. New header row: B1:E1 = years, F1 = "mean" (kept as text, matching
#        the source data's text-typed header labels). G1:I1 no longer used.
$ws.Range("B1").Value = "'2021"
$ws.Range("C1").Value = "'2022"
$ws.Range("D1").Value = "'2023"
$ws.Range("E1").Value = "'2024"
$ws.Range("F1").Value = "mean"
$ws.Range("G1:I1").Delete(-4159)  # xlShiftToLeft (scoped to row 1 only)

# --- 4. New row labels in column A (one row per statistic).
$ws.Range("A2").Value = "TB"
$ws.Range("A3").Value = "PC"
$ws.Range("A4").Value = "dS0"
$ws.Range("A5").Value = "dSF"
$ws.Range("A6").Value = "K"
$ws.Range("A7").Value = "IP"
$ws.Range("A8").Value = "I0"
$ws.Range("A9").Value = "IF"

# --- 5. Regenerated numeric data: rows = statistics, columns = years/mean.
# TB
$ws.Range("B2").Value = 2.06323212995964
$ws.Range("C2").Value = 2.174726767300435
$ws.Range("D2").Value = 2.242432058206409
$ws.Range("E2").Value = 2.112243226672423
$ws.Range("F2").Value = 2.148158545534727
# PC
$ws.Range("B3").Value = 24.17040945351953
$ws.Range("C3").Value = 24.68191017424561
$ws.Range("D3").Value = 24.79242123692412
$ws.Range("E3").Value = 23.9192487359343
$ws.Range("F3").Value = 24.39099740015589
# dS0
$ws.Range("B4").Value = 0.4225483524692952
$ws.Range("C4").Value = 0.4922814279462545
$ws.Range("D4").Value = 0.5867662140440464
$ws.Range("E4").Value = 0.3684276001125894
$ws.Range("F4").Value = 0.4675058986430464
# dSF
$ws.Range("B5").Value = 0.3262897972843836
$ws.Range("C5").Value = 0.2917839861045567
$ws.Range("D5").Value = 0.318710233571351
$ws.Range("E5").Value = 0.159799895900133
$ws.Range("F5").Value = 0.2741459782151061
# K
$ws.Range("B6").Value = 1.657087683535477
$ws.Range("C6").Value = 1.556289121901398
$ws.Range("D6").Value = 1.542387543798052
$ws.Range("E6").Value = 1.448562904744733
$ws.Range("F6").Value = 1.551081813494915
# IP
$ws.Range("B7").Value = 1.782514426289791
$ws.Range("C7").Value = 1.842098776136301
$ws.Range("D7").Value = 1.805920879557073
$ws.Range("E7").Value = 1.796504409142666
$ws.Range("F7").Value = 1.806759622781458
# I0
$ws.Range("B8").Value = 8.217485573710208
$ws.Range("C8").Value = 8.157901223863698
$ws.Range("D8").Value = 8.194079120442927
$ws.Range("E8").Value = 8.203495590857333
$ws.Range("F8").Value = 8.19324037721854
# IF
$ws.Range("B9").Value = 9
$ws.Range("C9").Value = 9
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = 9
$ws.Range("F9").Value = 9

# --- 6. Restore the original selection/view state (A1 active) now that the
#        new table occupies A1:F9.
$ws.Range("A1").Select()
